# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Row 16 now holds the single "TIBISAY MARIA PEREZ MORENO" (CC 50641270)
# record for period 1610 (Salario Basico stays 689455), while rows 17-33 are
# rewritten to hold "ALONSO POSADA SIERRA" (CC 73199752), one row per period
# from 1701 through 1805 (ascending order), with the Salario Basico
# (column G) updated from 689455 to 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: new worker record added ahead of Alonso Posada Sierra's rows.
$ws.Cells.Item(16, 3).Value = "50641270"
$ws.Cells.Item(16, 4).Value = "TIBISAY MARIA PEREZ MORENO"
$ws.Cells.Item(16, 5).Value = "1610"

$periods = @("1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805")

$row = 17
foreach ($periodo in $periods) {
    $ws.Cells.Item($row, 3).Value = "73199752"
    $ws.Cells.Item($row, 4).Value = "ALONSO POSADA SIERRA"
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 7).Value = 781242
    $row = $row + 1
}
